$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (B1:E1) - shared string values
$ws.Range("B1").Value = "layerB"
$ws.Range("C1").Value = "minicolumn"
$ws.Range("D1").Value = "hypercolumn"
$ws.Range("E1").Value = "layerA"

# Data rows 2-37
$ws.Range("A2").Value = 0.14851510907452498
$ws.Range("B2").Value = 5
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 4
$ws.Range("A3").Value = 0.14912278089866118
$ws.Range("B3").Value = 4
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 23
$ws.Range("A4").Value = 0.19037965090090636
$ws.Range("B4").Value = 5
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 23
$ws.Range("A5").Value = 0.1143382949222053
$ws.Range("B5").Value = 5
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 4
$ws.Range("A6").Value = 0.03784856577422081
$ws.Range("B6").Value = 4
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 23
$ws.Range("A7").Value = 0.1659653838406194
$ws.Range("B7").Value = 5
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 23
$ws.Range("A8").Value = 0.046328936639688435
$ws.Range("B8").Value = 5
$ws.Range("C8").Value = 2
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 4
$ws.Range("A9").Value = 0.04155309047786392
$ws.Range("B9").Value = 4
$ws.Range("C9").Value = 2
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 23
$ws.Range("A10").Value = 0.08681094419864815
$ws.Range("B10").Value = 5
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 23
$ws.Range("A11").Value = 0.12971370597969578
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = 3
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 4
$ws.Range("A12").Value = 0.04068361270046923
$ws.Range("B12").Value = 4
$ws.Range("C12").Value = 3
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 23
$ws.Range("A13").Value = 0.24919639115256959
$ws.Range("B13").Value = 5
$ws.Range("C13").Value = 3
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 23
$ws.Range("A14").Value = 0.06410720499822353
$ws.Range("B14").Value = 5
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 4
$ws.Range("A15").Value = 0.08489057724821769
$ws.Range("B15").Value = 4
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 23
$ws.Range("A16").Value = 0.18184077197824708
$ws.Range("B16").Value = 5
$ws.Range("C16").Value = 0
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 23
$ws.Range("A17").Value = 0.13022119317931488
$ws.Range("B17").Value = 5
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 4
$ws.Range("A18").Value = 0.03457546063131148
$ws.Range("B18").Value = 4
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 23
$ws.Range("A19").Value = 0.0592547579755059
$ws.Range("B19").Value = 5
$ws.Range("C19").Value = 1
$ws.Range("D19").Value = 1
$ws.Range("E19").Value = 23
$ws.Range("A20").Value = 0.060400417357177216
$ws.Range("B20").Value = 5
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 4
$ws.Range("A21").Value = 0.08597586519514165
$ws.Range("B21").Value = 4
$ws.Range("C21").Value = 2
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = 23
$ws.Range("A22").Value = 0.17908304986516244
$ws.Range("B22").Value = 5
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 23
$ws.Range("A23").Value = 0.20966374631556323
$ws.Range("B23").Value = 5
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 4
$ws.Range("A24").Value = 0.1477151116091055
$ws.Range("B24").Value = 4
$ws.Range("C24").Value = 3
$ws.Range("D24").Value = 1
$ws.Range("E24").Value = 23
$ws.Range("A25").Value = 0.19138675383684176
$ws.Range("B25").Value = 5
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 23
$ws.Range("A26").Value = 0.19121638516839604
$ws.Range("B26").Value = 5
$ws.Range("C26").Value = 0
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = 4
$ws.Range("A27").Value = 0.19215024408228326
$ws.Range("B27").Value = 4
$ws.Range("C27").Value = 0
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = 23
$ws.Range("A28").Value = 0.22143711641808975
$ws.Range("B28").Value = 5
$ws.Range("C28").Value = 0
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = 23
$ws.Range("A29").Value = 0.11683845688364769
$ws.Range("B29").Value = 5
$ws.Range("C29").Value = 1
$ws.Range("D29").Value = 2
$ws.Range("E29").Value = 4
$ws.Range("A30").Value = 0.2110846868384675
$ws.Range("B30").Value = 4
$ws.Range("C30").Value = 1
$ws.Range("D30").Value = 2
$ws.Range("E30").Value = 23
$ws.Range("A31").Value = 0.23635579864455308
$ws.Range("B31").Value = 5
$ws.Range("C31").Value = 1
$ws.Range("D31").Value = 2
$ws.Range("E31").Value = 23
$ws.Range("A32").Value = 0.11635337949191546
$ws.Range("B32").Value = 5
$ws.Range("C32").Value = 2
$ws.Range("D32").Value = 2
$ws.Range("E32").Value = 4
$ws.Range("A33").Value = 0.050765735646263285
$ws.Range("B33").Value = 4
$ws.Range("C33").Value = 2
$ws.Range("D33").Value = 2
$ws.Range("E33").Value = 23
$ws.Range("A34").Value = 0.16401741052103208
$ws.Range("B34").Value = 5
$ws.Range("C34").Value = 2
$ws.Range("D34").Value = 2
$ws.Range("E34").Value = 23
$ws.Range("A35").Value = 0.1709028814756272
$ws.Range("B35").Value = 5
$ws.Range("C35").Value = 3
$ws.Range("D35").Value = 2
$ws.Range("E35").Value = 4
$ws.Range("A36").Value = 0.13001694414860537
$ws.Range("B36").Value = 4
$ws.Range("C36").Value = 3
$ws.Range("D36").Value = 2
$ws.Range("E36").Value = 23
$ws.Range("A37").Value = 0.14854446591931197
$ws.Range("B37").Value = 5
$ws.Range("C37").Value = 3
$ws.Range("D37").Value = 2
$ws.Range("E37").Value = 23
